$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.321.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = "'1.665.44"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').Value = "'1.011"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.88%  '
$ws.Range('D5').Value = "'219.19"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = "'0.5350"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('D8').Value = "'0.2664"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('D9').Value = "'0.06406"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('D10').Value = "'20.76"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.87%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = "'4.566"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = "'1.664.46"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = "'1.893.83"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').Value = "'0.5536"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = "'0.0₅8211"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').Value = "'4.684"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').Value = "'10.30"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').Value = "'6.043"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').Value = "'146.52"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = "'7.210"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  +0.54%  '
$ws.Range('D28').Value = "'1.501"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.09%  '
$ws.Range('D29').Value = "'0.05850"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').Value = "'1.285"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = "'3.643"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.95%  '
$ws.Range('D32').Value = "'3.280"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = "'1.614"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('D34').Value = "'0.9685"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('D35').Value = "'2.828"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('D36').Value = "'2.418"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').Value = "'0.5824"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = "'0.8699"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('D40').Value = "'5.862"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('D41').Value = "'105.32"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').Value = "'1.052.85"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('D44').Value = "'1.804.64"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').Value = "'57.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').Value = "'1.015"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').Value = "'0.4389"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('D48').Value = "'7.997"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = "'0.0₈102"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.31%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.05170"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  -3.53%  '
